$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "/r/EarthPorn/comments/hhxc58/oc_amazing_sunset_colours_and_swell_in_south/"
$ws.Range("A4").Value = "/r/learnpython/comments/hhpl9c/ask_anything_monday_weekly_thread/"
$ws.Range("A5").Value = "/r/goodnews/comments/gwtp43/whats_new_content_creators_june_2020/"
$ws.Range("A6").Value = "/r/goodnews/comments/hg7813/good_news_its_friday_whats_your_feelgood_story/"
$ws.Range("A7").Value = "/r/learnpython/comments/hhu0zj/the_best_moment_ever_the_moment_when_it_finally/"
$ws.Range("A9").Value = "/r/pics/comments/hhydw1/oc45_years_ago_i_quit_my_stressful_desk_job_now_i/"
$ws.Range("A11").Value = "/r/Python/comments/hefa1t/whats_everyone_working_on_this_week/"

$ws.Range("E16").Select()
